$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (new Price text, or $null if unchanged) and new Volume(1h) text.
# NumNeedsText indicates the Price text looks like a plain number (e.g. "211.52")
# and must be forced to store as text (matching the source data), otherwise Excel
# would auto-convert it into a numeric value.
$updates = @(
    @{ Row = 2;  D = "26.656.06"; DIsNum = $false; E = "  -0.14%  " },
    @{ Row = 3;  D = "1.597.55";  DIsNum = $false; E = "  -0.01%  " },
    @{ Row = 4;  D = $null;       DIsNum = $false; E = "  +0.18%  " },
    @{ Row = 5;  D = "211.52";    DIsNum = $true;  E = "  +0.01%  " },
    @{ Row = 6;  D = $null;       DIsNum = $false; E = "  +0.65%  " },
    @{ Row = 7;  D = $null;       DIsNum = $false; E = "  +0.16%  " },
    @{ Row = 8;  D = $null;       DIsNum = $false; E = "  -0.06%  " },
    @{ Row = 9;  D = $null;       DIsNum = $false; E = "  +0.55%  " },
    @{ Row = 10; D = "19.57";     DIsNum = $true;  E = "  -0.90%  " },
    @{ Row = 11; D = $null;       DIsNum = $false; E = "  -0.11%  " },
    @{ Row = 12; D = "1.820.72";  DIsNum = $false; E = "  -0.04%  " },
    @{ Row = 13; D = "1.595.06";  DIsNum = $false; E = "  +0.02%  " },
    @{ Row = 14; D = "4.03";      DIsNum = $true;  E = "  -0.06%  " },
    @{ Row = 15; D = $null;       DIsNum = $false; E = "  +0.07%  " },
    @{ Row = 16; D = "65.16";     DIsNum = $true;  E = "  +0.09%  " },
    @{ Row = 17; D = "26.628.64"; DIsNum = $false; E = "  -0.26%  " },
    @{ Row = 18; D = "0.0₃0738";  DIsNum = $false; E = "  +1.32%  " },
    @{ Row = 19; D = "209.98";    DIsNum = $true;  E = "  -0.03%  " },
    @{ Row = 20; D = $null;       DIsNum = $false; E = "  +0.12%  " },
    @{ Row = 21; D = "7.01";      DIsNum = $true;  E = "  +3.97%  " },
    @{ Row = 22; D = $null;       DIsNum = $false; E = "  +0.45%  " },
    @{ Row = 23; D = $null;       DIsNum = $false; E = "  +1.11%  " },
    @{ Row = 24; D = "9.00";      DIsNum = $true;  E = "  +0.85%  " },
    @{ Row = 25; D = "145.17";    DIsNum = $true;  E = "  -1.00%  " },
    @{ Row = 26; D = $null;       DIsNum = $false; E = "  +0.08%  " },
    @{ Row = 27; D = $null;       DIsNum = $false; E = "  -1.01%  " },
    @{ Row = 28; D = $null;       DIsNum = $false; E = "  -0.67%  " },
    @{ Row = 29; D = "15.30";     DIsNum = $true;  E = "  -0.12%  " },
    @{ Row = 30; D = "0.0515";    DIsNum = $true;  E = "  +2.35%  " },
    @{ Row = 31; D = $null;       DIsNum = $false; E = "  +0.24%  " },
    @{ Row = 32; D = $null;       DIsNum = $false; E = "  +0.90%  " },
    @{ Row = 33; D = $null;       DIsNum = $false; E = "  +1.31%  " },
    @{ Row = 34; D = "1.282.65";  DIsNum = $false; E = "  -1.21%  " },
    @{ Row = 35; D = $null;       DIsNum = $false; E = "  -6.57%  " },
    @{ Row = 36; D = $null;       DIsNum = $false; E = "  +0.64%  " },
    @{ Row = 37; D = $null;       DIsNum = $false; E = "  +1.00%  " },
    @{ Row = 38; D = $null;       DIsNum = $false; E = "  -0.87%  " },
    @{ Row = 39; D = "0.835";     DIsNum = $true;  E = "  -0.98%  " },
    @{ Row = 40; D = $null;       DIsNum = $false; E = "  +19.38%  " },
    @{ Row = 41; D = $null;       DIsNum = $false; E = "  +2.27%  " },
    @{ Row = 42; D = $null;       DIsNum = $false; E = "  -0.14%  " },
    @{ Row = 43; D = "0.785";     DIsNum = $true;  E = "  -0.67%  " },
    @{ Row = 44; D = "63.71";     DIsNum = $true;  E = "  -0.06%  " },
    @{ Row = 45; D = "1.734.47";  DIsNum = $false; E = "  +0.02%  " },
    @{ Row = 46; D = "90.70";     DIsNum = $true;  E = "  +0.55%  " },
    @{ Row = 47; D = "1.58";      DIsNum = $true;  E = "  -3.10%  " },
    @{ Row = 48; D = $null;       DIsNum = $false; E = "  +1.45%  " },
    @{ Row = 49; D = "0.0508";    DIsNum = $true;  E = "  +0.62%  " },
    @{ Row = 50; D = $null;       DIsNum = $false; E = "  -0.17%  " },
    @{ Row = 51; D = "7.41";      DIsNum = $true;  E = "  -0.71%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($u.DIsNum) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
